# "Add files via upload" — re-upload of the evaluation sheet that swaps the
# evidence links from GitHub to Trello (the team moved their board), plus a
# tidy-up of the "Sprint 3" section header so it matches the plain style used
# by the other section headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Generales")

# --- Evidence column (D): "Video en Github" / "Screen en Github" -> Trello ---
$videoCells  = @("D6", "D7", "D8", "D15", "D16", "D17", "D18", "D24", "D25", "D26", "D27", "D28")
foreach ($addr in $videoCells) {
    $ws.Range($addr).Value = "Video en Trello"
}
$ws.Range("D9").Value = "Screen en Trello"

# --- "Sprint 3" header (A22) loses its stray highlighted style, matching A6/A13 ---
$ws.Range("A22").Font.Name = "Arial"
$ws.Range("A22").Font.Bold = $false
$ws.Range("A22").Font.Size = 11
$ws.Range("A22").HorizontalAlignment = -4131
$ws.Range("A22").VerticalAlignment = -4108
$ws.Range("A22").WrapText = $true

# --- Last worked range the author had selected before saving ---
[void]$ws.Range("D24:D28").Select()
